$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: KMeansClustering (replaces former DBSCAN row)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "KMeansClustering"
$ws.Range("C2").Value = 0.76907130982131477

# Row 3: MeanShiftClustering (new row)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "MeanShiftClustering"
$ws.Range("C3").Value = 0.76689637454639914

# Row 4: DBSCAN_Clustering (moved down, now last)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "DBSCAN_Clustering"
$ws.Range("C4").Value = 0.1440231944232267

# Copy the bordered/centered index-column formatting from A2 onto the
# newly-populated A3/A4 cells so the whole index column matches.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# Re-fit the columns now that the longer model names ("MeanShiftClustering",
# "KMeansClustering") are in place.
$ws.Columns.Item(1).ColumnWidth = 1.1
$ws.Columns.Item(2).ColumnWidth = 18.6
$ws.Columns.Item(3).ColumnWidth = 20.6

$ws.Range("A1").Select() | Out-Null
